# Scheduled-runner update: refresh cached market-board price/profit figures
# (currentAveragePrice*, Leve price & profit columns) on a handful of rows
# across several job sheets. Some rows also gain/lose a profit cell
# (M/N) as the underlying source data toggles between "has HQ data" and
# "NQ-only" states.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1907.7778
$ws.Range("I43").Value = 1845.125
$ws.Range("J43").Value = 1998.909
$ws.Range("K43").Value = 1845.125
$ws.Range("L43").Value = 1998.909
$ws.Range("M43").Value = -1776.125
$ws.Range("N43").Value = -2136.909

$ws.Range("H137").Value = 1697.2916
$ws.Range("I137").Value = 1346.75
$ws.Range("J137").Value = 3450
$ws.Range("K137").Value = 4040.25
$ws.Range("L137").Value = 10350
$ws.Range("M137").Value = -1490.25
$ws.Range("N137").Value = -15450

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1680.75
$ws.Range("I61").Value = 1142.0555
$ws.Range("J61").Value = 2650.4
$ws.Range("K61").Value = 1142.0555
$ws.Range("L61").Value = 2650.4
$ws.Range("M61").Value = -930.0554999999999
$ws.Range("N61").Value = -3074.4

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws.Range("H136").Value = 1680.75
$ws.Range("I136").Value = 1142.0555
$ws.Range("J136").Value = 2650.4
$ws.Range("K136").Value = 3426.1665
$ws.Range("L136").Value = 7951.200000000001
$ws.Range("M136").Value = -876.1664999999998
$ws.Range("N136").Value = -13051.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H43").Value = 270800.4
$ws.Range("J43").Value = 270800.4
$ws.Range("L43").Value = 270800.4
$ws.Range("N43").Value = -271162.4

$ws.Range("H94").Value = 688.05
$ws.Range("I94").Value = 505.72726
$ws.Range("J94").Value = 910.8889
$ws.Range("K94").Value = 505.72726
$ws.Range("L94").Value = 910.8889
$ws.Range("M94").Value = -54.72726
$ws.Range("N94").Value = -1812.8889

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 425.25
$ws.Range("I7").Value = 100
$ws.Range("J7").Value = 750.5
$ws.Range("K7").Value = 100
$ws.Range("L7").Value = 750.5
$ws.Range("M7").Value = 13
$ws.Range("N7").Value = -976.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 6046.909
$ws.Range("I3").Value = 5203.5
$ws.Range("J3").Value = 8296
$ws.Range("K3").Value = 15610.5
$ws.Range("L3").Value = 24888
$ws.Range("M3").Value = -15498.5
$ws.Range("N3").Value = -25112

$ws.Range("H17").Value = 633.1667
$ws.Range("I17").Value = 174.25
$ws.Range("J17").Value = 1551
$ws.Range("K17").Value = 522.75
$ws.Range("L17").Value = 4653
$ws.Range("M17").Value = -353.75
$ws.Range("N17").Value = -4991

$ws.Range("H64").Value = 2409.0908
$ws.Range("J64").Value = 2444.4443
$ws.Range("L64").Value = 7333.3329
$ws.Range("N64").Value = -7873.3329

$ws.Range("H67").Value = 2409.0908
$ws.Range("J67").Value = 2444.4443
$ws.Range("L67").Value = 7333.3329
$ws.Range("N67").Value = -9205.332900000001

$ws.Range("H69").Value = 2750
$ws.Range("I69").Value = 2000
$ws.Range("K69").Value = 6000
$ws.Range("M69").Value = -5189

$ws.Range("H70").Value = 102231.1
$ws.Range("I70").Value = 251577.75
$ws.Range("J70").Value = 2666.6667
$ws.Range("K70").Value = 754733.25
$ws.Range("L70").Value = 8000.000100000001
$ws.Range("M70").Value = -754418.25
$ws.Range("N70").Value = -8630.000100000001

$ws.Range("H72").Value = 2750
$ws.Range("I72").Value = 2000
$ws.Range("K72").Value = 18000
$ws.Range("M72").Value = -13944

$ws.Range("H73").Value = 102231.1
$ws.Range("I73").Value = 251577.75
$ws.Range("J73").Value = 2666.6667
$ws.Range("K73").Value = 754733.25
$ws.Range("L73").Value = 8000.000100000001
$ws.Range("M73").Value = -753641.25
$ws.Range("N73").Value = -10184.0001

$ws.Range("H80").Value = 34667.668
$ws.Range("J80").Value = 34667.668
$ws.Range("L80").Value = 104003.004
$ws.Range("N80").Value = -105875.004

$ws.Range("H83").Value = 34667.668
$ws.Range("J83").Value = 34667.668
$ws.Range("L83").Value = 312009.012
$ws.Range("N83").Value = -321369.012

$ws.Range("H114").Value = 424
$ws.Range("I114").Value = 424
$ws.Range("J114").Value = 0
$ws.Range("K114").Value = 1272
$ws.Range("L114").Value = 0
$ws.Range("M114").Value = 1982
$ws.Range("N114").ClearContents()

$ws.Range("H133").Value = 6350
$ws.Range("J133").Value = 12000
$ws.Range("L133").Value = 36000
$ws.Range("N133").Value = -46120

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 188.54546
$ws.Range("I2").Value = 152.66667
$ws.Range("J2").Value = 350
$ws.Range("K2").Value = 152.66667
$ws.Range("L2").Value = 350
$ws.Range("M2").Value = -39.66667000000001
$ws.Range("N2").Value = -576

$ws.Range("H35").Value = 10995
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 10995
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 10995
$ws.Range("N35").Value = -11591
$ws.Range("M35").ClearContents()

$ws.Range("H70").Value = 228301.89
$ws.Range("I70").Value = 337951.34
$ws.Range("J70").Value = 9003
$ws.Range("K70").Value = 337951.34
$ws.Range("L70").Value = 9003
$ws.Range("M70").Value = -337681.34
$ws.Range("N70").Value = -9543

$ws.Range("H73").Value = 228301.89
$ws.Range("I73").Value = 337951.34
$ws.Range("J73").Value = 9003
$ws.Range("K73").Value = 337951.34
$ws.Range("L73").Value = 9003
$ws.Range("M73").Value = -337015.34
$ws.Range("N73").Value = -10875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 7412512
$ws.Range("I16").Value = 9692808
$ws.Range("J16").Value = 1550.5
$ws.Range("K16").Value = 9692808
$ws.Range("L16").Value = 1550.5
$ws.Range("M16").Value = -9692638
$ws.Range("N16").Value = -1890.5

$ws.Range("H22").Value = 1017.4091
$ws.Range("I22").Value = 939.8
$ws.Range("K22").Value = 939.8
$ws.Range("M22").Value = -644.8

$ws.Range("H27").Value = 1017.4091
$ws.Range("I27").Value = 939.8
$ws.Range("K27").Value = 939.8
$ws.Range("M27").Value = -832.8

$ws.Range("H132").Value = 2289.9219
$ws.Range("I132").Value = 2259.44
$ws.Range("J132").Value = 2398.7856
$ws.Range("K132").Value = 6778.32
$ws.Range("L132").Value = 7196.3568
$ws.Range("M132").Value = -4248.32
$ws.Range("N132").Value = -12256.3568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 20550
$ws.Range("J92").Value = 20550
$ws.Range("L92").Value = 20550
$ws.Range("N92").Value = -25542
